$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 381, shifting the existing rows (old 381-408) down to 383-410.
$ws.Range("A381:R382").Insert()

# Copy formatting (esp. the date style in column D) from the row that is now 383 (formerly 381)
# into the two newly inserted blank rows, so the new rows look consistent with the rest of the table.
$ws.Range("A383:R384").Copy()
$ws.Range("A381:R382").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Row 381: new "Primera" quality record
$ws.Cells.Item(381, 1).Value = 8
$ws.Cells.Item(381, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(381, 3).Value = "Coquimbo"
$ws.Cells.Item(381, 4).Value = 44706
$ws.Cells.Item(381, 5).Value = 4
$ws.Cells.Item(381, 6).Value = 100112017
$ws.Cells.Item(381, 7).Value = "Apio"
$ws.Cells.Item(381, 8).Value = "Americana (o)"
$ws.Cells.Item(381, 9).Value = "Primera"
$ws.Cells.Item(381, 10).Value = 2460
$ws.Cells.Item(381, 11).Value = 8000
$ws.Cells.Item(381, 12).Value = 9000
$ws.Cells.Item(381, 13).Value = 8500
$ws.Cells.Item(381, 14).Value = "`$/docena de matas"
$ws.Cells.Item(381, 15).Value = "Provincia del Elqu" + [char]0x00ED
$ws.Cells.Item(381, 16).Value = 1417
$ws.Cells.Item(381, 17).Value = 6
$ws.Cells.Item(381, 18).Value = "Hortaliza"

# Row 382: new "Segunda" quality record
$ws.Cells.Item(382, 1).Value = 8
$ws.Cells.Item(382, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(382, 3).Value = "Coquimbo"
$ws.Cells.Item(382, 4).Value = 44706
$ws.Cells.Item(382, 5).Value = 4
$ws.Cells.Item(382, 6).Value = 100112017
$ws.Cells.Item(382, 7).Value = "Apio"
$ws.Cells.Item(382, 8).Value = "Americana (o)"
$ws.Cells.Item(382, 9).Value = "Segunda"
$ws.Cells.Item(382, 10).Value = 1320
$ws.Cells.Item(382, 11).Value = 6000
$ws.Cells.Item(382, 12).Value = 7000
$ws.Cells.Item(382, 13).Value = 6500
$ws.Cells.Item(382, 14).Value = "`$/docena de matas"
$ws.Cells.Item(382, 15).Value = "Provincia del Elqu" + [char]0x00ED
$ws.Cells.Item(382, 16).Value = 1083
$ws.Cells.Item(382, 17).Value = 6
$ws.Cells.Item(382, 18).Value = "Hortaliza"
